$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" (Strike count derived replacement) values for column G, rows 2-24
$values = @{
    2  = 5
    3  = 6
    4  = 7
    5  = 11
    6  = 11
    7  = 9
    8  = 5
    9  = 5
    10 = 8
    11 = 6
    12 = 4
    13 = 7
    14 = 9
    15 = 5
    16 = 6
    17 = 8
    18 = 7
    19 = 8
    20 = 4
    21 = 3
    22 = 3
    23 = 2
    24 = 4
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
